$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style from the existing header cell (H1) to the new header cells
# so they match the bold/bordered/centered header formatting.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Per-row data for columns I (I0) and J (IF)
$data = @(
    @(2,10,10),
    @(3,6,6),
    @(4,5,5),
    @(5,7,7),
    @(6,7,7),
    @(7,7,7),
    @(8,8,9),
    @(9,9,9),
    @(10,7,7),
    @(11,7,7),
    @(12,8,8),
    @(13,8,8),
    @(14,8,8),
    @(15,8,8),
    @(16,8,8),
    @(17,8,8),
    @(18,7,7),
    @(19,8,8),
    @(20,6,6),
    @(21,7,8),
    @(22,5,6),
    @(23,8,8),
    @(24,7,8),
    @(25,6,7),
    @(26,5,6),
    @(27,8,8),
    @(28,8,9),
    @(29,6,7),
    @(30,8,8),
    @(31,7,7),
    @(32,6,6),
    @(33,8,8),
    @(34,8,8),
    @(35,7,7),
    @(36,7,7),
    @(37,7,8),
    @(38,8,8),
    @(39,8,8),
    @(40,6,6),
    @(41,10,10),
    @(42,7,7),
    @(43,10,10),
    @(44,7,8),
    @(45,8,9),
    @(46,8,8),
    @(47,1,1),
    @(48,8,8),
    @(49,4,4),
    @(50,7,8),
    @(51,6,7),
    @(52,5,5),
    @(53,5,5)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $i0 = $entry[1]
    $iF = $entry[2]
    $ws.Cells.Item($r, 9).Value = $i0
    $ws.Cells.Item($r, 10).Value = $iF
}
